# Add two new columns "I0" (I) and "IF" (J) to the sheet, matching the
# existing header/style pattern and populating the data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style/format from the last existing header cell (H1)
# onto the two new header cells, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data values for rows 2-17, columns I (I0) and J (IF).
$iValues = @(9, 5, 8, 5, 8, 9, 6, 6, 1, 7, 7, 6, 8, 9, 8, 6)
$jValues = @(9, 6, 8, 5, 8, 9, 7, 7, 1, 7, 7, 6, 8, 9, 8, 6)

for ($r = 2; $r -le 17; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iValues[$idx]
    $ws.Cells.Item($r, 10).Value = $jValues[$idx]
}
